# Natmi following Dr Hou advice
# Update ligand/receptor expressing-cell counts and the resulting
# derived specificity / edge-weight metrics for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.240179
$ws.Range("H2").Value = 3.720537
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.940402333333334
$ws.Range("N2").Value = 23.821207
$ws.Range("O2").Value = 0.1931648990487216
$ws.Range("P2").Value = 0.1931648990487216
$ws.Range("Q2").Value = 9.847520225350999
$ws.Range("R2").Value = 88.627682028159
$ws.Range("S2").Value = 0.1931648990487216
$ws.Range("T2").Value = 0.1931648990487216

# --- Row 3 ---
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.240179
$ws.Range("H3").Value = 3.720537
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 26.95384733333333
$ws.Range("N3").Value = 80.861542
$ws.Range("O3").Value = 0.6557019380820612
$ws.Range("P3").Value = 0.6557019380820612
$ws.Range("Q3").Value = 33.42759543200599
$ws.Range("R3").Value = 300.848358888054
$ws.Range("S3").Value = 0.6557019380820612
$ws.Range("T3").Value = 0.6557019380820612

# --- Row 4 ---
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.240179
$ws.Range("H4").Value = 3.720537
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.212609666666666
$ws.Range("N4").Value = 18.637829
$ws.Range("O4").Value = 0.1511331628692172
$ws.Range("P4").Value = 0.1511331628692172
$ws.Range("Q4").Value = 7.704748043796999
$ws.Range("R4").Value = 69.342732394173
$ws.Range("S4").Value = 0.1511331628692172
$ws.Range("T4").Value = 0.1511331628692172
